# Burndown Chart update: fill in sprint data (dates, sprint id, totals, daily effort)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint Burndown Chart")

# --- Header block (row 1 / row 7): Date, Week, Group, Current Day ---
$ws.Range("I1").Value  = 44536      # Date
$ws.Range("K1").Value  = 3          # Week
$ws.Range("M1").Value  = 27         # Group
$ws.Range("I7").Value  = 44536      # Current Day (date)
$ws.Range("F7").Value  = 1

# --- Old decorative title row (row 4) is cleared out; row 5 becomes taller ---
$ws.Rows("4:4").Clear()
$ws.Rows("4:4").AutoFit()
$ws.Rows("5:5").RowHeight = 22.5

# --- Sprint backlog totals ---
$ws.Range("F8").Value = 60          # Total Estimated Sprint Backlog

# --- Daily "Actual Work" burned (column B), rows 12-41 ---
$actual = @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,5,10,4,3,4,4,10,8,5)
for ($i = 0; $i -lt $actual.Length; $i++) {
    $r = 12 + $i
    $ws.Cells.Item($r, 2).Value = $actual[$i]
}

# C12 is the seed (literal) of the cumulative "Done" column; C13:C41 stay formulas
$ws.Range("C12").Value = 0

# D12/D13 get re-entered (same formula, kept in sync with the rest of the column)
$ws.Range("D12").Formula = '=$C$41-C12'
$ws.Range("D13").Formula = '=$C$41-C13'

$wb.Application.Calculate()

# --- Final selection, as left by the author ---
$ws.Range("P28").Select() | Out-Null
